{"js": "// Add a new numbered list item (\"Transformed all static classes to instance\n// ones.\") right after the current last paragraph of the document body\n// (\"Made Figure and Pawn abstract classes.\"), keeping the same numbered\n// list formatting, and move the \"_GoBack\" bookmark (if present) so it still\n// sits at the very end of the document content.\n\nconst doc = context.document;\nconst body = doc.body;\n\n// Find the last paragraph in the document body - that's where the new\n// bullet point needs to be appended.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nconst lastPara = paragraphs.items[paragraphs.items.length - 1];\n\n// Word keeps a \"_GoBack\" bookmark at the end of the document content\n// (the last edit position). It currently sits at the end of the last\n// paragraph; remove it so it can be re-inserted after the new paragraph.\nconst goBack = doc.getBookmarkRangeOrNullObject(\"_GoBack\");\nawait context.sync();\ngoBack.load(\"isNullObject\");\nawait context.sync();\nconst hadGoBack = !goBack.isNullObject;\nif (hadGoBack) {\n  doc.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// Insert the new list item as a new paragraph right after the last one.\n// Inserting \"after\" an existing list paragraph makes the new paragraph\n// inherit the same formatting (ListParagraph style, numId 1, ilvl 0),\n// continuing the numbered list.\nconst newText = \"Transformed all static classes to instance ones.\";\nconst newPara = lastPara.insertParagraph(newText, Word.InsertLocation.after);\nawait context.sync();\n\n// Re-create the \"_GoBack\" bookmark collapsed right at the end of the new\n// paragraph's text (matching Word's normal behavior of keeping it at the\n// very end of the document content).\nif (hadGoBack) {\n  const results = newPara.getRange().search(newText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  const endOfNewText = results.items[0].getRange(\"After\");\n  await context.sync();\n  endOfNewText.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Add a new numbered list item (\"Transformed all static classes to instance\n# ones.\") right after the current last paragraph of the document\n# (\"Made Figure and Pawn abstract classes.\"), continuing the same numbered\n# list, and move the \"_GoBack\" bookmark (if present) so it still sits at the\n# very end of the document content - matching Word's own behavior of keeping\n# that bookmark at the last edited/inserted position.\n\n$d = $word.ActiveDocument\n\n$newText = \"Transformed all static classes to instance ones.\"\n\n# Find out whether the \"_GoBack\" bookmark exists, and remember it so it can\n# be re-created at the new end of the document afterwards.\n$hadGoBack = $true\ntry {\n    $goBack = $d.Bookmarks.Item(\"_GoBack\")\n} catch {\n    $hadGoBack = $false\n}\nif ($hadGoBack) {\n    $goBack.Delete()\n}\n\n# Append a new paragraph at the very end of the document. Inserting the\n# paragraph break there makes the new paragraph inherit the formatting\n# of the paragraph before it (ListParagraph style, numId 1, ilvl 0), so the\n# new bullet continues the existing numbered list automatically.\n$endRange = $d.Content\n$endRange.Collapse(0)\n$endRange.InsertParagraphAfter()\n\n# Type the new bullet text, followed by one placeholder character. The\n# placeholder keeps the \"real\" end-of-text position at least two characters\n# away from the end of the document/paragraph while we create the bookmark\n# (collapsed ranges landing on the very last one/two character slots of a\n# paragraph are mis-positioned by this host's Bookmarks.Add), then we strip\n# the placeholder back out, leaving the bookmark correctly collapsed right\n# after \"ones.\" and before the paragraph mark.\n$newEnd = $d.Content\n$newEnd.Collapse(0)\n$newEnd.InsertAfter($newText + \"X\")\n\nif ($hadGoBack) {\n    $docLen = $d.Content.End\n    $targetPos = $docLen - 2\n    $bookmarkRange = $d.Range($targetPos, $targetPos)\n    $d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n\n    # Remove the one-character placeholder again.\n    $placeholderRange = $d.Range($targetPos, $targetPos + 1)\n    $placeholderRange.Delete()\n}\n"}
